$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2024-05-13 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-14 Tuesday", 2) | Out-Null

# Update the 100 arithmetic answers in the 20x5 table (row-major order)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "41-24=17"
$t.Cell(1, 2).Range.Text = "44+16=60"
$t.Cell(1, 3).Range.Text = "11+9=20"
$t.Cell(1, 4).Range.Text = "74-6=68"
$t.Cell(1, 5).Range.Text = "5+33=38"

$t.Cell(2, 1).Range.Text = "15+67=82"
$t.Cell(2, 2).Range.Text = "16+60=76"
$t.Cell(2, 3).Range.Text = "73-60=13"
$t.Cell(2, 4).Range.Text = "26+62=88"
$t.Cell(2, 5).Range.Text = "53-35=18"

$t.Cell(3, 1).Range.Text = "62-1=61"
$t.Cell(3, 2).Range.Text = "13+29=42"
$t.Cell(3, 3).Range.Text = "32+12=44"
$t.Cell(3, 4).Range.Text = "51-10=41"
$t.Cell(3, 5).Range.Text = "34-24=10"

$t.Cell(4, 1).Range.Text = "9+60=69"
$t.Cell(4, 2).Range.Text = "75-32=43"
$t.Cell(4, 3).Range.Text = "21-14=7"
$t.Cell(4, 4).Range.Text = "4+15=19"
$t.Cell(4, 5).Range.Text = "90-41=49"

$t.Cell(5, 1).Range.Text = "65-47=18"
$t.Cell(5, 2).Range.Text = "35+63=98"
$t.Cell(5, 3).Range.Text = "72+25=97"
$t.Cell(5, 4).Range.Text = "49-17=32"
$t.Cell(5, 5).Range.Text = "87-44=43"

$t.Cell(6, 1).Range.Text = "69+25=94"
$t.Cell(6, 2).Range.Text = "42+32=74"
$t.Cell(6, 3).Range.Text = "42+30=72"
$t.Cell(6, 4).Range.Text = "10+34=44"
$t.Cell(6, 5).Range.Text = "21+78=99"

$t.Cell(7, 1).Range.Text = "91-77=14"
$t.Cell(7, 2).Range.Text = "23-5=18"
$t.Cell(7, 3).Range.Text = "73-3=70"
$t.Cell(7, 4).Range.Text = "89-83=6"
$t.Cell(7, 5).Range.Text = "36+10=46"

$t.Cell(8, 1).Range.Text = "68-52=16"
$t.Cell(8, 2).Range.Text = "63-60=3"
$t.Cell(8, 3).Range.Text = "91-32=59"
$t.Cell(8, 4).Range.Text = "33+20=53"
$t.Cell(8, 5).Range.Text = "98-89=9"

$t.Cell(9, 1).Range.Text = "3+45=48"
$t.Cell(9, 2).Range.Text = "92-25=67"
$t.Cell(9, 3).Range.Text = "87-23=64"
$t.Cell(9, 4).Range.Text = "2+61=63"
$t.Cell(9, 5).Range.Text = "71+26=97"

$t.Cell(10, 1).Range.Text = "3+52=55"
$t.Cell(10, 2).Range.Text = "56+43=99"
$t.Cell(10, 3).Range.Text = "99-92=7"
$t.Cell(10, 4).Range.Text = "51-30=21"
$t.Cell(10, 5).Range.Text = "68-33=35"

$t.Cell(11, 1).Range.Text = "49-8=41"
$t.Cell(11, 2).Range.Text = "10+18=28"
$t.Cell(11, 3).Range.Text = "14-12=2"
$t.Cell(11, 4).Range.Text = "24-19=5"
$t.Cell(11, 5).Range.Text = "27+3=30"

$t.Cell(12, 1).Range.Text = "68+7=75"
$t.Cell(12, 2).Range.Text = "55-13=42"
$t.Cell(12, 3).Range.Text = "5+28=33"
$t.Cell(12, 4).Range.Text = "38-28=10"
$t.Cell(12, 5).Range.Text = "72-11=61"

$t.Cell(13, 1).Range.Text = "35-12=23"
$t.Cell(13, 2).Range.Text = "52+32=84"
$t.Cell(13, 3).Range.Text = "78-32=46"
$t.Cell(13, 4).Range.Text = "47+52=99"
$t.Cell(13, 5).Range.Text = "45-25=20"

$t.Cell(14, 1).Range.Text = "74+4=78"
$t.Cell(14, 2).Range.Text = "26+42=68"
$t.Cell(14, 3).Range.Text = "2+40=42"
$t.Cell(14, 4).Range.Text = "93-53=40"
$t.Cell(14, 5).Range.Text = "65+16=81"

$t.Cell(15, 1).Range.Text = "5+34=39"
$t.Cell(15, 2).Range.Text = "80+3=83"
$t.Cell(15, 3).Range.Text = "57-44=13"
$t.Cell(15, 4).Range.Text = "6+57=63"
$t.Cell(15, 5).Range.Text = "38+61=99"

$t.Cell(16, 1).Range.Text = "71+17=88"
$t.Cell(16, 2).Range.Text = "0+59=59"
$t.Cell(16, 3).Range.Text = "99-13=86"
$t.Cell(16, 4).Range.Text = "26+39=65"
$t.Cell(16, 5).Range.Text = "52+24=76"

$t.Cell(17, 1).Range.Text = "1+78=79"
$t.Cell(17, 2).Range.Text = "64-5=59"
$t.Cell(17, 3).Range.Text = "12-8=4"
$t.Cell(17, 4).Range.Text = "61-33=28"
$t.Cell(17, 5).Range.Text = "98-12=86"

$t.Cell(18, 1).Range.Text = "12+38=50"
$t.Cell(18, 2).Range.Text = "23-19=4"
$t.Cell(18, 3).Range.Text = "43+31=74"
$t.Cell(18, 4).Range.Text = "38-30=8"
$t.Cell(18, 5).Range.Text = "44+21=65"

$t.Cell(19, 1).Range.Text = "11+87=98"
$t.Cell(19, 2).Range.Text = "32-25=7"
$t.Cell(19, 3).Range.Text = "76+20=96"
$t.Cell(19, 4).Range.Text = "38+57=95"
$t.Cell(19, 5).Range.Text = "91-90=1"

$t.Cell(20, 1).Range.Text = "34-12=22"
$t.Cell(20, 2).Range.Text = "67+31=98"
$t.Cell(20, 3).Range.Text = "33+30=63"
$t.Cell(20, 4).Range.Text = "78-71=7"
$t.Cell(20, 5).Range.Text = "95-10=85"
